# Update gh-pages to output generated at 456a3b4
# Refresh the "want to go" counts (column F) that bilibili reports for each
# event, flag the cancelled AF convention, and mark its ticket column as
# "not for sale" instead of a price — applied to the "展览" sheet, the
# "演出" sheet, and their aggregate in "全部类型".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览" (Exhibitions)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("展览")

$ws.Range("F2").Value  = 20655
$ws.Range("F5").Value  = 1112
$ws.Range("F7").Value  = 7728
$ws.Range("F10").Value = 291
$ws.Range("F12").Value = 176
$ws.Range("F13").Value = 141

$ws.Range("C15").Value = "苏州·AF特摄/动漫/游戏三厨狂欢节（取消）"
$ws.Range("G15").Value = "不可售"

$ws.Range("F16").Value = 212
$ws.Range("F18").Value = 481
$ws.Range("F19").Value = 77
$ws.Range("F22").Value = 75
$ws.Range("F25").Value = 1153
$ws.Range("F28").Value = 199
$ws.Range("F29").Value = 5214
$ws.Range("F31").Value = 102
$ws.Range("F32").Value = 4924
$ws.Range("F36").Value = 12859
$ws.Range("F37").Value = 1349
$ws.Range("F38").Value = 104
$ws.Range("F40").Value = 63
$ws.Range("F41").Value = 291
$ws.Range("F42").Value = 401

# ---------------------------------------------------------------------
# Sheet "演出" (Performances)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")

$ws2.Range("F2").Value = 244
$ws2.Range("F4").Value = 36

# ---------------------------------------------------------------------
# Sheet "全部类型" (All types combined)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F2").Value  = 20655
$ws4.Range("F5").Value  = 1112
$ws4.Range("F7").Value  = 7728
$ws4.Range("F10").Value = 291
$ws4.Range("F12").Value = 176
$ws4.Range("F13").Value = 141

$ws4.Range("C15").Value = "苏州·AF特摄/动漫/游戏三厨狂欢节（取消）"
$ws4.Range("G15").Value = "不可售"

$ws4.Range("F16").Value = 212
$ws4.Range("F18").Value = 481
$ws4.Range("F19").Value = 77
$ws4.Range("F22").Value = 75
$ws4.Range("F25").Value = 1153
$ws4.Range("F28").Value = 199
$ws4.Range("F29").Value = 244
$ws4.Range("F30").Value = 5214
$ws4.Range("F33").Value = 102
$ws4.Range("F34").Value = 36
$ws4.Range("F35").Value = 4924
$ws4.Range("F39").Value = 12859
$ws4.Range("F40").Value = 1349
$ws4.Range("F41").Value = 104
$ws4.Range("F43").Value = 63
$ws4.Range("F44").Value = 291
$ws4.Range("F45").Value = 401
